$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M: "Objeto" (the bid's object/description), pulled from existing
# data for a couple of rows. Values are entered bottom-up then the header,
# matching how the shared-string table ends up ordered.
$ws.Range("M5").Value = "Abalone"
$ws.Range("M4").Value = "Jiribatuba2"
$ws.Range("M2").Value = "objeto teste"
$ws.Range("M1").Value = "Objeto"

# Column M width.
$ws.Columns.Item(13).ColumnWidth = 24.125

# Touch formatting on M4 so a dedicated (plain) style record is written for it.
$ws.Range("M4").WrapText = $False

# Move the selection to the newly-populated M2 cell.
$null = $ws.Range("M2").Select()

# Page setup (paper size + orientation) for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
